# Apply updated cryptocurrency price/volume figures (and the row-41/42
# TerraClassic <-> FTXToken swap) to match the refreshed data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.577.91"
$ws.Range("E2").Value = "  +4.64%  "
$ws.Range("D3").Value = "2.220.34"
$ws.Range("E3").Value = "  +2.80%  "
$ws.Range("E4").Value = "  +0.13%  "
$ws.Range("D5").Value = "'228.05"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("D6").Value = "'0.620"
$ws.Range("E6").Value = "  -1.37%  "
$ws.Range("D7").Value = "'60.86"
$ws.Range("E7").Value = "  -4.07%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("D9").Value = "'0.400"
$ws.Range("E9").Value = "  +2.00%  "
$ws.Range("D10").Value = "'58.09"
$ws.Range("E10").Value = "  -0.06%  "
$ws.Range("E11").Value = "  +3.21%  "
$ws.Range("E12").Value = "  -0.21%  "
$ws.Range("D13").Value = "2.554.15"
$ws.Range("E13").Value = "  +3.00%  "
$ws.Range("D14").Value = "'15.62"
$ws.Range("E14").Value = "  -1.98%  "
$ws.Range("D15").Value = "'21.45"
$ws.Range("E15").Value = "  -2.02%  "
$ws.Range("D16").Value = "'0.793"
$ws.Range("E16").Value = "  -1.42%  "
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").Value = "2.224.27"
$ws.Range("E18").Value = "  +2.89%  "
$ws.Range("D19").Value = "41.486.66"
$ws.Range("E19").Value = "  +4.66%  "
$ws.Range("D20").Value = "'72.41"
$ws.Range("E20").Value = "  +0.90%  "
$ws.Range("D21").Value = "0.0₃0882"
$ws.Range("E21").Value = "  +4.44%  "
$ws.Range("D22").Value = "'6.01"
$ws.Range("E22").Value = "  -1.02%  "
$ws.Range("D23").Value = "'246.16"
$ws.Range("E23").Value = "  +6.78%  "
$ws.Range("E24").Value = "  -0.05%  "
$ws.Range("D25").Value = "'2.35"
$ws.Range("E25").Value = "  -1.01%  "
$ws.Range("D26").Value = "'2.27"
$ws.Range("E26").Value = "  -3.93%  "
$ws.Range("D27").Value = "'9.54"
$ws.Range("E27").Value = "  -0.10%  "
$ws.Range("D28").Value = "'167.90"
$ws.Range("E28").Value = "  -2.53%  "
$ws.Range("D29").Value = "'0.139"
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("D30").Value = "'19.83"
$ws.Range("E30").Value = "  -0.35%  "
$ws.Range("D31").Value = "'1.41"
$ws.Range("E31").Value = "  -4.18%  "
$ws.Range("D32").Value = "'2.64"
$ws.Range("E32").Value = "  -1.57%  "
$ws.Range("D33").Value = "'0.121"
$ws.Range("E33").Value = "  -0.72%  "
$ws.Range("E34").Value = "  +7.89%  "
$ws.Range("D35").Value = "'4.62"
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D36").Value = "'0.0620"
$ws.Range("E36").Value = "  +0.11%  "
$ws.Range("D37").Value = "'6.58"
$ws.Range("E37").Value = "  -5.28%  "
$ws.Range("D38").Value = "'3.68"
$ws.Range("E38").Value = "  +2.22%  "
$ws.Range("E39").Value = "  -1.58%  "
$ws.Range("D40").Value = "'0.998"
$ws.Range("E40").Value = "  -0.21%  "
$ws.Range("B41").Value = "FTXToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D41").Value = "'4.84"
$ws.Range("E41").Value = "  -3.83%  "
$ws.Range("B42").Value = "TerraClassic"
$ws.Range("C42").Value = "https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc"
$ws.Range("D42").Value = "'0.000233"
$ws.Range("E42").Value = "  +25.99%  "
$ws.Range("E43").Value = "  +13.40%  "
$ws.Range("E44").Value = "  +3.83%  "
$ws.Range("D45").Value = "'0.0975"
$ws.Range("E45").Value = "  +6.09%  "
$ws.Range("D46").Value = "'98.72"
$ws.Range("E46").Value = "  -4.25%  "
$ws.Range("D47").Value = "1.465.99"
$ws.Range("E47").Value = "  -3.21%  "
$ws.Range("E48").Value = "  -3.12%  "
$ws.Range("D49").Value = "'2.78"
$ws.Range("E49").Value = "  -0.83%  "
$ws.Range("D50").Value = "'16.28"
$ws.Range("E50").Value = "  -7.20%  "
$ws.Range("D51").Value = "'1.07"
$ws.Range("E51").Value = "  -1.75%  "
